$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: bold formatting (new font + cell style) ---
$ws.Range("B1:H1").Font.Bold = $true

# --- Row 19 used to read "solpreneur"; a new "deep learning" row is
#     inserted right after it and the old row is relabelled
#     "solopreneur". Write the two new labels in this order so the
#     shared-string table ends up with "deep learning" before
#     "solopreneur" (matching the reused/new string slots), then fix
#     up which row points at which label. ---
$ws.Range("A19").Value = "deep learning"
$ws.Range("A20").Value = "solopreneur"
$ws.Range("A19").Value = "solopreneur"
$ws.Range("A20").Value = "deep learning"

$ws.Range("B20").Value = 3
$ws.Range("C20").Value = 4
$ws.Range("F20").Value = 5

# --- Append new survey rows 21-24 ---
$ws.Range("A21").Value = "upskill"
$ws.Range("B21").Value = 5
$ws.Range("C21").Value = 4
$ws.Range("F21").Value = 5

$ws.Range("A22").Value = "deepfake"
$ws.Range("C22").Value = 4
$ws.Range("F22").Value = 5

$ws.Range("A23").Value = "followership"
$ws.Range("B23").Value = 4
$ws.Range("C23").Value = 4
$ws.Range("F23").Value = 5

$ws.Range("A24").Value = "Brexiter / Brexiteer"
$ws.Range("B24").Value = 5
$ws.Range("F24").Value = 5

# --- Freeze header row, scroll down, and reselect like the author left it ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A10").Select()
$ws.Range("D25").Select()
